$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 307

# Columns (by index) that may contain HYPERLINK formulas whose display text
# needs to be added as a second argument.
$linkCols = @(19, 20, 22, 23, 24, 25)   # S, T, V, W, X, Y

for ($row = $firstRow; $row -le $lastRow; $row++) {

    # Column C ("Förändrad" / last-changed date) moves from 45184 to 45186
    # for every data row.
    $cCell = $ws.Cells.Item($row, 3)
    $cCell.Value = 45186

    # The display label used inside HYPERLINK(...) is the row's
    # "Beteckning" value stored in column A.
    $label = $ws.Cells.Item($row, 1).Value()

    foreach ($col in $linkCols) {
        $cell = $ws.Cells.Item($row, $col)
        $formula = $cell.Formula()
        if ($formula.StartsWith("=HYPERLINK(") -and -not $formula.Contains(",")) {
            $newFormula = $formula.Substring(0, $formula.Length - 1) + ', "' + $label + '")'
            $cell.Formula = $newFormula
        }
    }
}
